# Generate Report for Handback
# Replace the two source file UUID-based names / associated xlf + timestamp
# values across the Overview / zh-cn / de-de sheets.
#
# NOTE: the hyperlink *targets* (external GitHub blob URLs, recorded in the
# worksheet .rels and referenced by the unchanged r:id values in the diff)
# keep the OLD source file names - only the visible cell text / hyperlink
# display text is renamed to the new source file names.

$wb = $excel.ActiveWorkbook

$oldName1 = "b1f45bce-ead6-4194-8875-91118ab5a876"
$oldName2 = "b2fd5750-ae99-46e5-8811-4c47ddc44e2f"
$newName1 = "3f31608e-0d5b-405d-b49e-acdd59fef6b8"
$newName2 = "ffffa091a720-58b1-49a1-a246-fdf430061668"

$newZhCnXlf = "$newName1.ae2eb6082684251e7dc10f7ff4183c6dc6d4e931.zh-cn.xlf"
$newDeDeXlf = "$newName1.ae2eb6082684251e7dc10f7ff4183c6dc6d4e931.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newName1.md"
$wsOverview.Range("B2").Value = "e2e\$newName1.md"
$wsOverview.Range("G2").Value = "2016-08-26 07:01:21"

$wsOverview.Range("A3").Value = "$newName2.md"
$wsOverview.Range("B3").Value = "e2e\$newName2.md"
$wsOverview.Range("G3").Value = "2016-08-26 07:01:21"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/$oldName1.md", "", "", "e2e\$newName1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/$oldName2.md", "", "", "e2e\$newName2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newName1.md"
$wsZhCn.Range("I2").Value = "$newName1.md"
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = "2016-08-26 07:01:16"
$wsZhCn.Range("J2").Value = $newZhCnXlf
$wsZhCn.Range("K2").Value = "2016-08-26 07:01:32"

$wsZhCn.Range("A3").Value = "$newName2.md"
$wsZhCn.Range("I3").Value = "$newName2.md"
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("H3").Value = "2016-08-26 07:01:16"
$wsZhCn.Range("J3").Value = $newZhCnXlf
$wsZhCn.Range("K3").Value = "2016-08-26 07:01:32"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/$oldName1.md", "", "", "$newName1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ef46ede68f2d74634ec11861dbcafbb3067099a0/e2e/$oldName1.md", "", "", "$newName1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/$oldName2.md", "", "", "$newName2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ef46ede68f2d74634ec11861dbcafbb3067099a0/e2e/$oldName2.md", "", "", "$newName2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newName1.md"
$wsDeDe.Range("I2").Value = "$newName1.md"
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = "2016-08-26 07:01:21"
$wsDeDe.Range("J2").Value = $newDeDeXlf
$wsDeDe.Range("K2").Value = "2016-08-26 07:01:39"

$wsDeDe.Range("A3").Value = "$newName2.md"
$wsDeDe.Range("I3").Value = "$newName2.md"
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("H3").Value = "2016-08-26 07:01:21"
$wsDeDe.Range("J3").Value = $newDeDeXlf
$wsDeDe.Range("K3").Value = "2016-08-26 07:01:39"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/$oldName1.md", "", "", "$newName1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/20e635be66b1ab09f2b0f5add6766d98bf7b57f5/e2e/$oldName1.md", "", "", "$newName1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2035fcaa69ebe218bc5d7533e92c0b294afaa35d/e2e/$oldName2.md", "", "", "$newName2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/20e635be66b1ab09f2b0f5add6766d98bf7b57f5/e2e/$oldName2.md", "", "", "$newName2.md")
